$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.223728813559322
$ws.Range("C2").Value = 0.5322033898305085
$ws.Range("J2").Value = 0.006779661016949152
$ws.Range("P2").Value = 0.1796610169491525
$ws.Range("S2").Value = 0.0576271186440678
$ws.Range("B3").Value = 0.00625
$ws.Range("C3").Value = 0.01875
$ws.Range("J3").Value = 0.025
$ws.Range("P3").Value = 0.7875
$ws.Range("S3").Value = 0.1625
$ws.Range("J4").Value = 0.0425531914893617
$ws.Range("P4").Value = 0.7446808510638298
$ws.Range("S4").Value = 0.2127659574468085
$ws.Range("B6").Value = 0.02463054187192118
$ws.Range("D6").Value = 0.004926108374384237
$ws.Range("F6").Value = 0.04433497536945813
$ws.Range("J6").Value = 0.1773399014778325
$ws.Range("O6").Value = 0.02955665024630542
$ws.Range("Q6").Value = 0.2315270935960591
$ws.Range("R6").Value = 0.103448275862069
$ws.Range("S6").Value = 0.3842364532019704
$ws.Range("B7").Value = 0.1373626373626374
$ws.Range("D7").Value = 0.005494505494505495
$ws.Range("E7").Value = 0.005494505494505495
$ws.Range("F7").Value = 0.08241758241758242
$ws.Range("J7").Value = 0.0989010989010989
$ws.Range("O7").Value = 0.01098901098901099
$ws.Range("Q7").Value = 0.2142857142857143
$ws.Range("R7").Value = 0.09340659340659341
$ws.Range("S7").Value = 0.3516483516483517
$ws.Range("B8").Value = 0.08470588235294117
$ws.Range("D8").Value = 0.0188235294117647
$ws.Range("F8").Value = 0.04470588235294118
$ws.Range("J8").Value = 0.1129411764705882
$ws.Range("O8").Value = 0.01411764705882353
$ws.Range("Q8").Value = 0.2
$ws.Range("R8").Value = 0.08705882352941176
$ws.Range("S8").Value = 0.4376470588235294
$ws.Range("B9").Value = 0.1396648044692737
$ws.Range("D9").Value = 0.0111731843575419
$ws.Range("F9").Value = 0.0670391061452514
$ws.Range("J9").Value = 0.106145251396648
$ws.Range("O9").Value = 0.0111731843575419
$ws.Range("Q9").Value = 0.1620111731843575
$ws.Range("R9").Value = 0.0893854748603352
$ws.Range("S9").Value = 0.4134078212290503
$ws.Range("B10").Value = 0.1095317725752508
$ws.Range("D10").Value = 0.02926421404682274
$ws.Range("F10").Value = 0.08528428093645485
$ws.Range("J10").Value = 0.132943143812709
$ws.Range("O10").Value = 0.00919732441471572
$ws.Range("Q10").Value = 0.2090301003344482
$ws.Range("R10").Value = 0.08444816053511706
$ws.Range("S10").Value = 0.3403010033444816
$ws.Range("G11").Value = 0.1490909090909091
$ws.Range("J11").Value = 0.07636363636363637
$ws.Range("K11").Value = 0.2036363636363636
$ws.Range("L11").Value = 0.5309090909090909
$ws.Range("S11").Value = 0.04
$ws.Range("G12").Value = 0.7397260273972602
$ws.Range("J12").Value = 0.1986301369863014
$ws.Range("L12").Value = 0.02054794520547945
$ws.Range("S12").Value = 0.0410958904109589
$ws.Range("G13").Value = 0.7058823529411765
$ws.Range("J13").Value = 0.2549019607843137
$ws.Range("S13").Value = 0.0392156862745098
$ws.Range("F15").Value = 0.01714285714285714
$ws.Range("H15").Value = 0.1542857142857143
$ws.Range("I15").Value = 0.06285714285714286
$ws.Range("J15").Value = 0.3771428571428572
$ws.Range("K15").Value = 0.08
$ws.Range("M15").Value = 0.01714285714285714
$ws.Range("O15").Value = 0.02285714285714286
$ws.Range("S15").Value = 0.2685714285714286
$ws.Range("F16").Value = 0.004950495049504951
$ws.Range("H16").Value = 0.1584158415841584
$ws.Range("I16").Value = 0.09405940594059406
$ws.Range("J16").Value = 0.405940594059406
$ws.Range("K16").Value = 0.1237623762376238
$ws.Range("M16").Value = 0.0198019801980198
$ws.Range("O16").Value = 0.05445544554455446
$ws.Range("S16").Value = 0.1386138613861386
$ws.Range("F17").Value = 0.01342281879194631
$ws.Range("H17").Value = 0.1923937360178971
$ws.Range("I17").Value = 0.1029082774049217
$ws.Range("J17").Value = 0.4138702460850112
$ws.Range("K17").Value = 0.06935123042505593
$ws.Range("M17").Value = 0.02237136465324385
$ws.Range("O17").Value = 0.04697986577181208
$ws.Range("S17").Value = 0.1387024608501119
$ws.Range("F18").Value = 0.005208333333333333
$ws.Range("H18").Value = 0.21875
$ws.Range("I18").Value = 0.1041666666666667
$ws.Range("J18").Value = 0.3541666666666667
$ws.Range("K18").Value = 0.07291666666666667
$ws.Range("M18").Value = 0.03125
$ws.Range("O18").Value = 0.07291666666666667
$ws.Range("S18").Value = 0.140625
$ws.Range("F19").Value = 0.005912162162162162
$ws.Range("H19").Value = 0.2043918918918919
$ws.Range("I19").Value = 0.07347972972972973
$ws.Range("J19").Value = 0.3817567567567567
$ws.Range("K19").Value = 0.1097972972972973
$ws.Range("M19").Value = 0.02364864864864865
$ws.Range("O19").Value = 0.06418918918918919
$ws.Range("S19").Value = 0.1368243243243243

Write-Output "Applied 106 cell updates"